$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 01:04"

# --- Update daily stats for a handful of existing countries (no name changes) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1320044
$ws.Range("C4").Value = 27421
$ws.Range("E4").Value = 1019507
$ws.Range("G4").Value = 1601
$ws.Range("H4").Value = 78529

# Ecuador (row 21)
$ws.Range("E21").Value = 25161
$ws.Range("G21").Value = 50
$ws.Range("H21").Value = 1704

# Chequia (row 49)
$ws.Range("B49").Value = 8077
$ws.Range("C49").Value = 46
$ws.Range("E49").Value = 3396
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 273

# Noruega (row 50)
$ws.Range("B50").Value = 8070
$ws.Range("C50").Value = 36
$ws.Range("E50").Value = 7820

# Venezuela (row 128)
$ws.Range("B128").Value = 388
$ws.Range("C128").Value = 7
$ws.Range("D128").Value = 190
$ws.Range("E128").Value = 188

# --- Countries block (rows 150-154): a new "Haiti" entry is added right before
# --- "Gibraltar", and the old "Haiti" row (which used to sit between "Brunei"
# --- and "Camboya") is removed. Net effect on this block of rows: each cell's
# --- country name and figures are updated in place to their new values.

# Row 150: now Haiti (was Gibraltar)
$ws.Range("A150").Value = "Haiti"
$ws.Range("B150").Value = 146
$ws.Range("C150").Value = 17
$ws.Range("D150").Value = 17
$ws.Range("E150").Value = 117
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 12

# Row 151: now Gibraltar (was Togo)
$ws.Range("A151").Value = "Gibraltar"
$ws.Range("B151").Value = 146
$ws.Range("C151").Value = 2
$ws.Range("D151").Value = 142
$ws.Range("E151").Value = 4
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

# Row 152: now Togo (was Guayana Francesa)
$ws.Range("A152").Value = "Togo"
$ws.Range("B152").Value = 141
$ws.Range("C152").Value = 10
$ws.Range("D152").Value = 85
$ws.Range("E152").Value = 50
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 10

# Row 153: now Guayana Francesa (was Brunei)
$ws.Range("A153").Value = "Guayana Francesa"
$ws.Range("B153").Value = 141
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 113
$ws.Range("E153").Value = 27
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 1

# Row 154: now Brunei (was Haiti)
$ws.Range("A154").Value = "Brunei"
$ws.Range("B154").Value = 141
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 132
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 2
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 1
